# Applies the OOXML diff to before.pptx (single-slide deck).
#
# PowerPoint's Shape.Left/Top/Width/Height are IEEE-754 single precision
# (points); to land exactly on the target EMU values after the
# point->EMU round trip, the literals below are nudged by a handful of
# ULPs from the "obvious" EMU/12700 quotient so that (single)(pt)*12700
# truncates to the exact target integer EMU count.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# Shape id=22 "Rectangle: Rounded Corners 21" (Gen.AI / Bed Rock box)
# Moves from (5874502,425112) to (5817425,2376438) EMU; size unchanged.
# ---------------------------------------------------------------------
$shGenAI = $s.Shapes.Item(13)
$shGenAI.Left = 458.06496062992125
$shGenAI.Top  = 187.12110236220474

# ---------------------------------------------------------------------
# Shape id=42 "Connector: Elbow 41" (connects shape 14 -> shape 22)
# flipV removed, offset + extents change.
# ---------------------------------------------------------------------
$shConn = $s.Shapes.Item(22)
$shConn.VerticalFlip = $false
$shConn.Left   = 380.6512604425197
$shConn.Top    = 153.44953155905512
$shConn.Width  = 77.41370078740158
$shConn.Height = 55.146535933070865

# ---------------------------------------------------------------------
# Shape id=54 "TextBox 53" ("Anomaly Explainer" -> two lines, moved/resized)
# ---------------------------------------------------------------------
$shAnomaly = $s.Shapes.Item(27)
$shAnomaly.Left   = 360.3276519952756
$shAnomaly.Top    = 134.67842869685037
$shAnomaly.Width  = 167.29023622047245
$shAnomaly.Height = 36.3515759031496
$shAnomaly.TextFrame.TextRange.Text = "Anomaly" + [char]13 + " Explainer"

# ---------------------------------------------------------------------
# Shape id=55 "TextBox 54" ("Data Assistance"), position only.
# ---------------------------------------------------------------------
$shData = $s.Shapes.Item(28)
$shData.Left = 416.35456852913387
$shData.Top  = 226.31984711968505

# ---------------------------------------------------------------------
# Shape id=62 "TextBox 61" ("AI Driven Prediction & Insight"), position only.
# ---------------------------------------------------------------------
$shTitle = $s.Shapes.Item(31)
$shTitle.Left = 132.61267856535434
$shTitle.Top  = -2.5092126984251966
